$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / metadata updates (rows 3,4,6,7) ---
# C3 ("Ficha de Caracterizacion") is stored as TEXT even though it looks numeric
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1262805"
$ws.Range("C4").Value = "CONTABILIDAD Y FINANZAS"
$ws.Range("C6").Value = 42639
$ws.Range("C7").Value = 43368

# --- Student table updates (rows 11-27 modified, 28-47 appended) ---
# Row 11
$ws.Range("A11").Value = "CC"
$ws.Range("B11").Value = 1005181647
$ws.Range("C11").Value = "JOAN MANUEL"
$ws.Range("D11").Value = "PALENCIA LOPEZ"
$ws.Range("E11").Value = "CERTIFICADO"

# Row 12
$ws.Range("A12").Value = "CC"
$ws.Range("B12").Value = 1020495471
$ws.Range("C12").Value = "JOAN MANUEL"
$ws.Range("D12").Value = "PALENCIA LOPEZ"
$ws.Range("E12").Value = "CERTIFICADO"

# Row 13
$ws.Range("A13").Value = "CC"
$ws.Range("B13").Value = 1095946925
$ws.Range("C13").Value = "DANNA GERALDINE"
$ws.Range("D13").Value = "VEGA SANCHEZ"
$ws.Range("E13").Value = "CERTIFICADO"

# Row 14
$ws.Range("A14").Value = "CC"
$ws.Range("B14").Value = 1096189854
$ws.Range("C14").Value = "LIZETH PAOLA"
$ws.Range("D14").Value = "LOBO SALCEDO"
$ws.Range("E14").Value = "TRASLADADO"

# Row 15
$ws.Range("A15").Value = "CC"
$ws.Range("B15").Value = 1096194400
$ws.Range("C15").Value = "KATHERINE"
$ws.Range("D15").Value = "TERRAZA ALVAREZ"
$ws.Range("E15").Value = "CERTIFICADO"

# Row 16
$ws.Range("A16").Value = "CC"
$ws.Range("B16").Value = 1096197940
$ws.Range("C16").Value = "MAYERLY PAOLA"
$ws.Range("D16").Value = "CAÑA PALACIO"
$ws.Range("E16").Value = "CERTIFICADO"

# Row 17
$ws.Range("A17").Value = "CC"
$ws.Range("B17").Value = 1096201993
$ws.Range("C17").Value = "JOHANNA PAOLA"
$ws.Range("D17").Value = "BARRERA RODRIGUEZ"
$ws.Range("E17").Value = "CERTIFICADO"

# Row 18
$ws.Range("A18").Value = "CC"
$ws.Range("B18").Value = 1096208435
$ws.Range("C18").Value = "KAROL VANESSA"
$ws.Range("D18").Value = "ORTEGA TOLOZA"
$ws.Range("E18").Value = "CERTIFICADO"

# Row 19
$ws.Range("A19").Value = "CC"
$ws.Range("B19").Value = 1096210258
$ws.Range("C19").Value = "KAREN MARCELA"
$ws.Range("D19").Value = "ALVARADO GARCIA"
$ws.Range("E19").Value = "CERTIFICADO"

# Row 20
$ws.Range("A20").Value = "CC"
$ws.Range("B20").Value = 1096211115
$ws.Range("C20").Value = "PAOLA ANDREA"
$ws.Range("D20").Value = "MERCADO MARIN"
$ws.Range("E20").Value = "CERTIFICADO"

# Row 21
$ws.Range("A21").Value = "CC"
$ws.Range("B21").Value = 1096213910
$ws.Range("C21").Value = "YENIFER PAOLA"
$ws.Range("D21").Value = "LEYVA ALVARINO"
$ws.Range("E21").Value = "CERTIFICADO"

# Row 22
$ws.Range("A22").Value = "CC"
$ws.Range("B22").Value = 1096214167
$ws.Range("C22").Value = "MARYI LISETH"
$ws.Range("D22").Value = "CHIMA TRIANA"
$ws.Range("E22").Value = "RETIRO VOLUNTARIO"

# Row 23
$ws.Range("A23").Value = "CC"
$ws.Range("B23").Value = 1096216042
$ws.Range("C23").Value = "CHERIL ANDREA"
$ws.Range("D23").Value = "NAVARRO GOMEZ"
$ws.Range("E23").Value = "CERTIFICADO"

# Row 24
$ws.Range("A24").Value = "CC"
$ws.Range("B24").Value = 1096223641
$ws.Range("C24").Value = "ELSA PAOLA"
$ws.Range("D24").Value = "VIDES OROZCO"
$ws.Range("E24").Value = "CERTIFICADO"

# Row 25
$ws.Range("A25").Value = "CC"
$ws.Range("B25").Value = 1096223976
$ws.Range("C25").Value = "KAROL"
$ws.Range("D25").Value = "BUITRAGO RINCON"
$ws.Range("E25").Value = "CERTIFICADO"

# Row 26
$ws.Range("A26").Value = "CC"
$ws.Range("B26").Value = 1096226289
$ws.Range("C26").Value = "FANNY PAOLA"
$ws.Range("D26").Value = "VALENCIA OSSES"
$ws.Range("E26").Value = "CERTIFICADO"

# Row 27
$ws.Range("A27").Value = "CC"
$ws.Range("B27").Value = 1096227259
$ws.Range("C27").Value = "JUAN DAVID"
$ws.Range("D27").Value = "JOYA BELLO"
$ws.Range("E27").Value = "CERTIFICADO"

# Row 28
$ws.Range("A28").Value = "CC"
$ws.Range("B28").Value = 1096229358
$ws.Range("C28").Value = "LITH JHAJAIRA"
$ws.Range("D28").Value = "PUERTA GOMEZ"
$ws.Range("E28").Value = "CANCELADO"

# Row 29
$ws.Range("A29").Value = "CC"
$ws.Range("B29").Value = 1096231912
$ws.Range("C29").Value = "LEIDY JOHANNA"
$ws.Range("D29").Value = "PEREIRA GARCIA"
$ws.Range("E29").Value = "CANCELADO"

# Row 30
$ws.Range("A30").Value = "CC"
$ws.Range("B30").Value = 1096233614
$ws.Range("C30").Value = "KATHERIN"
$ws.Range("D30").Value = "DIAZ CASTILLO"
$ws.Range("E30").Value = "CERTIFICADO"

# Row 31
$ws.Range("A31").Value = "CC"
$ws.Range("B31").Value = 1096234226
$ws.Range("C31").Value = "ANDRES FELIPE"
$ws.Range("D31").Value = "PEREZ SARMIENTO"
$ws.Range("E31").Value = "CANCELADO"

# Row 32
$ws.Range("A32").Value = "CC"
$ws.Range("B32").Value = 1096237824
$ws.Range("C32").Value = "JHON MAYRO"
$ws.Range("D32").Value = "PATIÑO CASTILLO"
$ws.Range("E32").Value = "TRASLADADO"

# Row 33
$ws.Range("A33").Value = "CC"
$ws.Range("B33").Value = 1096238898
$ws.Range("C33").Value = "MANUEL YAIR"
$ws.Range("D33").Value = "SILVA DURAN"
$ws.Range("E33").Value = "CERTIFICADO"

# Row 34
$ws.Range("A34").Value = "CC"
$ws.Range("B34").Value = 1096240013
$ws.Range("C34").Value = "INGRID VANESSA"
$ws.Range("D34").Value = "ARIAS DIAZ"
$ws.Range("E34").Value = "CANCELADO"

# Row 35
$ws.Range("A35").Value = "CC"
$ws.Range("B35").Value = 1096242025
$ws.Range("C35").Value = "YURLEY TATIANA"
$ws.Range("D35").Value = "CARDONA GARZON"
$ws.Range("E35").Value = "CERTIFICADO"

# Row 36
$ws.Range("A36").Value = "CC"
$ws.Range("B36").Value = 1096242447
$ws.Range("C36").Value = "EYLEEN YARITZA"
$ws.Range("D36").Value = "GARCIA RUEDA"
$ws.Range("E36").Value = "CERTIFICADO"

# Row 37
$ws.Range("A37").Value = "CC"
$ws.Range("B37").Value = 1096248773
$ws.Range("C37").Value = "WINDRY LISNETH"
$ws.Range("D37").Value = "ARENAS COLMENARES"
$ws.Range("E37").Value = "CERTIFICADO"

# Row 38
$ws.Range("A38").Value = "CC"
$ws.Range("B38").Value = 1096249526
$ws.Range("C38").Value = "JESSICA MELISSA"
$ws.Range("D38").Value = "CAMPO MONCADA"
$ws.Range("E38").Value = "CERTIFICADO"

# Row 39
$ws.Range("A39").Value = "CC"
$ws.Range("B39").Value = 1096252843
$ws.Range("C39").Value = "ABEL RICARDO"
$ws.Range("D39").Value = "MARTINEZ RUEDA"
$ws.Range("E39").Value = "CERTIFICADO"

# Row 40
$ws.Range("A40").Value = "CC"
$ws.Range("B40").Value = 1098607019
$ws.Range("C40").Value = "SIRLEY"
$ws.Range("D40").Value = "DELGADILLO SIERRA"
$ws.Range("E40").Value = "TRASLADADO"

# Row 41
$ws.Range("A41").Value = "CC"
$ws.Range("B41").Value = 1098695412
$ws.Range("C41").Value = "RODRIGO ANDRES"
$ws.Range("D41").Value = "PEREZ CARRASCAL"
$ws.Range("E41").Value = "CERTIFICADO"

# Row 42
$ws.Range("A42").Value = "CC"
$ws.Range("B42").Value = 1098769158
$ws.Range("C42").Value = "BRAYAN SNEIDER"
$ws.Range("D42").Value = "GAMARRA ZAPATA"
$ws.Range("E42").Value = "CERTIFICADO"

# Row 43
$ws.Range("A43").Value = "CC"
$ws.Range("B43").Value = 28020924
$ws.Range("C43").Value = "YAZMIH LORENA"
$ws.Range("D43").Value = "MONTES GALVAN"
$ws.Range("E43").Value = "CANCELADO"

# Row 44
$ws.Range("A44").Value = "CC"
$ws.Range("B44").Value = 37577019
$ws.Range("C44").Value = "SANDRA PATRICIA"
$ws.Range("D44").Value = "RAMIREZ ECHEVERRY"
$ws.Range("E44").Value = "CERTIFICADO"

# Row 45
$ws.Range("A45").Value = "CC"
$ws.Range("B45").Value = 37578408
$ws.Range("C45").Value = "YEIMI"
$ws.Range("D45").Value = "AMAYA SUAREZ"
$ws.Range("E45").Value = "CANCELADO"

# Row 46
$ws.Range("A46").Value = "CC"
$ws.Range("B46").Value = 37580170
$ws.Range("C46").Value = "VIVIANA"
$ws.Range("D46").Value = "RODRIGUEZ MUÑOZ"
$ws.Range("E46").Value = "CERTIFICADO"

# Row 47
$ws.Range("A47").Value = "CC"
$ws.Range("B47").Value = 63472623
$ws.Range("C47").Value = "YADIRA"
$ws.Range("D47").Value = "HOSTIA SALAS"
$ws.Range("E47").Value = "CERTIFICADO"

Write-Host "Applied novedades update"
